# Applies the cryptos-list price/volume refresh described by the commit diff.
# Cell values are plain text in the workbook (inlineStr), so numeric-looking
# "Price" values are written with a leading apostrophe to force Excel to keep
# them as text instead of auto-converting to a number (mirrors typing them in
# the Excel UI with a text/quote prefix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.178.81"
$ws.Range("E2").Value = "  -1.19%  "

$ws.Range("D3").Value = "1.894.20"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'245.83"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").Value = "'0.686"
$ws.Range("E6").Value = "  +8.27%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "'40.44"
$ws.Range("E8").Value = "  -4.55%  "

$ws.Range("E9").Value = "  +1.96%  "

$ws.Range("D10").Value = "'53.14"
$ws.Range("E10").Value = "  +11.32%  "

$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("D12").Value = "'0.0983"
$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D13").Value = "2.169.70"
$ws.Range("E13").Value = "  -0.82%  "

$ws.Range("D14").Value = "'12.57"
$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("D15").Value = "'0.702"
$ws.Range("E15").Value = "  +1.36%  "

$ws.Range("D16").Value = "1.902.40"
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").Value = "'4.77"
$ws.Range("E17").Value = "  -1.56%  "

$ws.Range("D18").Value = "35.179.39"
$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("D19").Value = "'71.87"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").Value = "0.0₃0814"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").Value = "'240.09"
$ws.Range("E21").Value = "  -1.51%  "

$ws.Range("D22").Value = "'12.52"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("D25").Value = "'2.30"
$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("D26").Value = "'2.36"
$ws.Range("E26").Value = "  +11.67%  "

$ws.Range("D27").Value = "'167.59"
$ws.Range("E27").Value = "  -2.27%  "

$ws.Range("D28").Value = "'8.52"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'18.21"
$ws.Range("E29").Value = "  +1.09%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.130"
$ws.Range("E30").Value = "  +3.36%  "

$ws.Range("D32").Value = "'4.14"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").Value = "'0.0563"
$ws.Range("E33").Value = "  -0.50%  "

$ws.Range("E34").Value = "  -0.22%  "

$ws.Range("D35").Value = "'1.87"
$ws.Range("E35").Value = "  +9.08%  "

$ws.Range("D36").Value = "'4.07"
$ws.Range("E36").Value = "  -2.04%  "

$ws.Range("D37").Value = "'0.903"
$ws.Range("E37").Value = "  -5.11%  "

$ws.Range("E38").Value = "  +11.76%  "

$ws.Range("D39").Value = "'2.01"
$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0208"
$ws.Range("E40").Value = "  +0.98%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.09"
$ws.Range("E41").Value = "  -2.62%  "

$ws.Range("D42").Value = "'0.0639"
$ws.Range("E42").Value = "  +7.52%  "

$ws.Range("D43").Value = "'16.09"
$ws.Range("E43").Value = "  +5.54%  "

$ws.Range("D44").Value = "'89.45"
$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("D45").Value = "1.342.31"
$ws.Range("E45").Value = "  -1.22%  "

$ws.Range("E46").Value = "  +2.29%  "

$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("D49").Value = "'45.21"
$ws.Range("E49").Value = "  -8.53%  "

$ws.Range("D50").Value = "'12.13"
$ws.Range("E50").Value = "  -7.23%  "

$ws.Range("D51").Value = "'6.42"
$ws.Range("E51").Value = "  -3.74%  "
